# Scheduled runner update: refresh currentAveragePrice/Leve price/profit
# columns (H-N) on the Ifrit_Profits job sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)
# with latest market-data values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3853.3635
$ws.Range("I32").Value = 4947.5
$ws.Range("J32").Value = 3228.1428
$ws.Range("K32").Value = 4947.5
$ws.Range("L32").Value = 3228.1428
$ws.Range("M32").Value = -4621.5
$ws.Range("N32").Value = -3880.1428
$ws.Range("H111").Value = 2925.5715
$ws.Range("I111").Value = 3519.75
$ws.Range("J111").Value = 2133.3333
$ws.Range("K111").Value = 10559.25
$ws.Range("L111").Value = 6399.999899999999
$ws.Range("M111").Value = -7492.25
$ws.Range("N111").Value = -12533.9999
$ws.Range("H127").Value = 1273.6842
$ws.Range("J127").Value = 1473.3334
$ws.Range("L127").Value = 4420.0002
$ws.Range("N127").Value = -14340.0002
$ws.Range("H129").Value = 835.88
$ws.Range("I129").Value = 518.2
$ws.Range("J129").Value = 915.3
$ws.Range("K129").Value = 1554.6
$ws.Range("L129").Value = 2745.9
$ws.Range("M129").Value = 3445.4
$ws.Range("N129").Value = -12745.9
$ws.Range("H132").Value = 3054.353
$ws.Range("I132").Value = 3131.7878
$ws.Range("J132").Value = 499
$ws.Range("K132").Value = 9395.3634
$ws.Range("L132").Value = 1497
$ws.Range("M132").Value = -6865.3634
$ws.Range("N132").Value = -6557
$ws.Range("H137").Value = 27030016
$ws.Range("I137").Value = 1799.4286
$ws.Range("J137").Value = 43481972
$ws.Range("K137").Value = 5398.2858
$ws.Range("L137").Value = 130445916
$ws.Range("M137").Value = -2848.2858
$ws.Range("N137").Value = -130451016
$ws.Range("H140").Value = 50646.668
$ws.Range("J140").Value = 50646.668
$ws.Range("L140").Value = 50646.668
$ws.Range("N140").Value = -61006.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 27782876
$ws.Range("I61").Value = 37040870
$ws.Range("J61").Value = 8890
$ws.Range("K61").Value = 37040870
$ws.Range("L61").Value = 8890
$ws.Range("M61").Value = -37040658
$ws.Range("N61").Value = -9314
$ws.Range("H74").Value = 18525512
$ws.Range("I74").Value = 26317116
$ws.Range("J74").Value = 20451.75
$ws.Range("K74").Value = 26317116
$ws.Range("L74").Value = 20451.75
$ws.Range("M74").Value = -26316242
$ws.Range("N74").Value = -22199.75
$ws.Range("H77").Value = 18525512
$ws.Range("I77").Value = 26317116
$ws.Range("J77").Value = 20451.75
$ws.Range("K77").Value = 131585580
$ws.Range("L77").Value = 102258.75
$ws.Range("M77").Value = -131581212
$ws.Range("N77").Value = -110994.75
$ws.Range("H110").Value = 1190.12
$ws.Range("I110").Value = 1041.7894
$ws.Range("K110").Value = 1041.7894
$ws.Range("M110").Value = 1003.2106
$ws.Range("H136").Value = 27782876
$ws.Range("I136").Value = 37040870
$ws.Range("J136").Value = 8890
$ws.Range("K136").Value = 111122610
$ws.Range("L136").Value = 26670
$ws.Range("M136").Value = -111120060
$ws.Range("N136").Value = -31770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1039
$ws.Range("I16").Value = 980.13336
$ws.Range("J16").Value = 1333.3334
$ws.Range("K16").Value = 980.13336
$ws.Range("L16").Value = 1333.3334
$ws.Range("M16").Value = -693.13336
$ws.Range("N16").Value = -1907.3334
$ws.Range("H58").Value = 3026.6167
$ws.Range("I58").Value = 1237.7222
$ws.Range("J58").Value = 3793.2856
$ws.Range("K58").Value = 1237.7222
$ws.Range("L58").Value = 3793.2856
$ws.Range("M58").Value = -1034.7222
$ws.Range("N58").Value = -4199.2856
$ws.Range("H113").Value = 1039
$ws.Range("I113").Value = 980.13336
$ws.Range("J113").Value = 1333.3334
$ws.Range("K113").Value = 980.13336
$ws.Range("L113").Value = 1333.3334
$ws.Range("M113").Value = 1189.86664
$ws.Range("N113").Value = -5673.3334
$ws.Range("H133").Value = 31560.857
$ws.Range("J133").Value = 31560.857
$ws.Range("L133").Value = 31560.857
$ws.Range("N133").Value = -36620.857
$ws.Range("H134").Value = 2296.3428
$ws.Range("I134").Value = 2208.9583
$ws.Range("K134").Value = 6626.874899999999
$ws.Range("M134").Value = -4091.874899999999
$ws.Range("H136").Value = 3026.6167
$ws.Range("I136").Value = 1237.7222
$ws.Range("J136").Value = 3793.2856
$ws.Range("K136").Value = 3713.1666
$ws.Range("L136").Value = 11379.8568
$ws.Range("M136").Value = -1163.1666
$ws.Range("N136").Value = -16479.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 1033.2858
$ws.Range("J45").Value = 1006.6
$ws.Range("L45").Value = 3019.8
$ws.Range("N45").Value = -4083.8
$ws.Range("H117").Value = 1101.9
$ws.Range("I117").Value = 866.1667
$ws.Range("J117").Value = 1455.5
$ws.Range("K117").Value = 2598.5001
$ws.Range("L117").Value = 4366.5
$ws.Range("M117").Value = 843.4998999999998
$ws.Range("N117").Value = -11250.5
$ws.Range("H129").Value = 1369.7273
$ws.Range("I129").Value = 622.7143
$ws.Range("J129").Value = 1718.3334
$ws.Range("K129").Value = 1868.1429
$ws.Range("L129").Value = 5155.0002
$ws.Range("M129").Value = 3131.8571
$ws.Range("N129").Value = -15155.0002
$ws.Range("H131").Value = 3432.7144
$ws.Range("J131").Value = 2103.6
$ws.Range("L131").Value = 6310.799999999999
$ws.Range("N131").Value = -16390.8
$ws.Range("H139").Value = 41668630
$ws.Range("I139").Value = 41668630
$ws.Range("K139").Value = 125005890
$ws.Range("M139").Value = -125000750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1730.32
$ws.Range("I126").Value = 1482.9231
$ws.Range("J126").Value = 1998.3334
$ws.Range("K126").Value = 4448.7693
$ws.Range("L126").Value = 5995.0002
$ws.Range("M126").Value = -1978.7693
$ws.Range("N126").Value = -10935.0002
$ws.Range("H132").Value = 2095.4688
$ws.Range("I132").Value = 1837.96
$ws.Range("J132").Value = 3015.1428
$ws.Range("K132").Value = 5513.88
$ws.Range("L132").Value = 9045.428400000001
$ws.Range("M132").Value = -2983.88
$ws.Range("N132").Value = -14105.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2379.35
$ws.Range("I40").Value = 2111.2
$ws.Range("J40").Value = 2647.5
$ws.Range("K40").Value = 2111.2
$ws.Range("L40").Value = 2647.5
$ws.Range("M40").Value = -1975.2
$ws.Range("N40").Value = -2919.5
$ws.Range("H68").Value = 2172
$ws.Range("I68").Value = 1900
$ws.Range("J68").Value = 2466.6667
$ws.Range("K68").Value = 1900
$ws.Range("L68").Value = 2466.6667
$ws.Range("M68").Value = -1151
$ws.Range("N68").Value = -3964.6667
$ws.Range("H71").Value = 2172
$ws.Range("I71").Value = 1900
$ws.Range("J71").Value = 2466.6667
$ws.Range("K71").Value = 9500
$ws.Range("L71").Value = 12333.3335
$ws.Range("M71").Value = -5756
$ws.Range("N71").Value = -19821.3335
$ws.Range("H82").Value = 1714.2858
$ws.Range("I82").Value = 1581.8182
$ws.Range("J82").Value = 2200
$ws.Range("K82").Value = 1581.8182
$ws.Range("L82").Value = 2200
$ws.Range("M82").Value = -1220.8182
$ws.Range("N82").Value = -2922
$ws.Range("H85").Value = 1714.2858
$ws.Range("I85").Value = 1581.8182
$ws.Range("J85").Value = 2200
$ws.Range("K85").Value = 1581.8182
$ws.Range("L85").Value = 2200
$ws.Range("M85").Value = -333.8181999999999
$ws.Range("N85").Value = -4696
$ws.Range("H136").Value = 2087.8125
$ws.Range("I136").Value = 1310.6
$ws.Range("J136").Value = 3383.1667
$ws.Range("K136").Value = 3931.8
$ws.Range("L136").Value = 10149.5001
$ws.Range("M136").Value = -1381.8
$ws.Range("N136").Value = -15249.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3995.5386
$ws.Range("I62").Value = 4483.6665
$ws.Range("J62").Value = 3577.1428
$ws.Range("K62").Value = 4483.6665
$ws.Range("L62").Value = 3577.1428
$ws.Range("M62").Value = -3859.6665
$ws.Range("N62").Value = -4825.1428
$ws.Range("H65").Value = 3995.5386
$ws.Range("I65").Value = 4483.6665
$ws.Range("J65").Value = 3577.1428
$ws.Range("K65").Value = 22418.3325
$ws.Range("L65").Value = 17885.714
$ws.Range("M65").Value = -19298.3325
$ws.Range("N65").Value = -24125.714
$ws.Range("H96").Value = 8335154.5
$ws.Range("I96").Value = 12501838
$ws.Range("J96").Value = 1787.5
$ws.Range("K96").Value = 12501838
$ws.Range("L96").Value = 1787.5
$ws.Range("M96").Value = -12500465
$ws.Range("N96").Value = -4533.5
$ws.Range("H136").Value = 7505
$ws.Range("I136").Value = 10176.786
$ws.Range("J136").Value = 2161.4285
$ws.Range("K136").Value = 30530.358
$ws.Range("L136").Value = 6484.2855
$ws.Range("M136").Value = -27980.358
$ws.Range("N136").Value = -11584.2855
